$d = $word.ActiveDocument
$ns  = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$rpr = '<w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="20"/></w:rPr>'

# ---------------------------------------------------------------------------
# Change 1: "Table names cannot have spaces ..." -> split into three runs:
#   "Table names" + " and column names" + " cannot have spaces ..."
# ---------------------------------------------------------------------------
$target = "Table names cannot have spaces and must be English-language letters. They also cannot be punctuation characters."
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13) -eq $target) {
        $start = $p.Range.Start
        $end   = $p.Range.End - 1
        $r     = $d.Range($start, $end)
        $run1  = "<w:r>$rpr<w:t>Table names</w:t></w:r>"
        $run2  = "<w:r>$rpr<w:t xml:space=`"preserve`"> and column names</w:t></w:r>"
        $run3  = "<w:r>$rpr<w:t xml:space=`"preserve`"> cannot have spaces and must be English-language letters. They also cannot be punctuation characters.</w:t></w:r>"
        $xml   = "<w:p $ns>" + $run1 + $run2 + $run3 + "</w:p>"
        $r.InsertXML($xml)
        break
    }
}

# ---------------------------------------------------------------------------
# Helper: append a brand-new ListParagraph (ilvl=2, numId=2) example line
# right after the paragraph whose text equals $afterText.
# $runsXml is the inner XML (runs / proofErr) for the new paragraph.
# ---------------------------------------------------------------------------
function Add-ExampleParagraph($afterText, $runsXml) {
    $d = $word.ActiveDocument
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.TrimEnd([char]13) -eq $afterText) {
            $start = $p.Range.Start
            $end   = $p.Range.End
            $full  = $d.Range($start, $end)

            # Recreate the (unchanged) original paragraph verbatim ...
            $origPPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:ind w:right="27"/>' + $rpr + '</w:pPr>'
            $origRun = "<w:r>$rpr<w:t>$afterText</w:t></w:r>"
            $origP   = "<w:p>" + $origPPr + $origRun + "</w:p>"

            # ... followed by the brand-new example paragraph.
            $newPPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr><w:ind w:right="27"/>' + $rpr + '</w:pPr>'
            $newP   = "<w:p>" + $newPPr + $runsXml + "</w:p>"

            $xml = "<w:body $ns>" + $origP + $newP + "</w:body>"
            $full.InsertXML($xml)
            break
        }
    }
}

# ---------------------------------------------------------------------------
# Change 2: after "Can delete all rows in a table" -> "DELETE ALL ROWS tablename"
# ---------------------------------------------------------------------------
$runs2 = "<w:r>$rpr<w:t xml:space=`"preserve`">DELETE ALL ROWS </w:t></w:r>" +
         '<w:proofErr w:type="spellStart"/>' +
         "<w:r>$rpr<w:t>tablename</w:t></w:r>" +
         '<w:proofErr w:type="spellEnd"/>'
Add-ExampleParagraph "Can delete all rows in a table" $runs2

# ---------------------------------------------------------------------------
# Change 3: after "Can delete an entire table" -> "DELETE TABLE tablename"
# ---------------------------------------------------------------------------
$runs3 = "<w:r>$rpr<w:t xml:space=`"preserve`">DELETE TABLE </w:t></w:r>" +
         '<w:proofErr w:type="spellStart"/>' +
         "<w:r>$rpr<w:t>tablename</w:t></w:r>" +
         '<w:proofErr w:type="spellEnd"/>'
Add-ExampleParagraph "Can delete an entire table" $runs3

# ---------------------------------------------------------------------------
# Change 4: after "The user must specify the row they wish to delete via the
# Primary Key" -> "DELETE primarykey FROM tablename"
# ---------------------------------------------------------------------------
$runs4 = "<w:r>$rpr<w:t xml:space=`"preserve`">DELETE </w:t></w:r>" +
         '<w:proofErr w:type="spellStart"/>' +
         "<w:r>$rpr<w:t>primarykey</w:t></w:r>" +
         '<w:proofErr w:type="spellEnd"/>' +
         "<w:r>$rpr<w:t xml:space=`"preserve`"> FROM tablename</w:t></w:r>"
Add-ExampleParagraph "The user must specify the row they wish to delete via the Primary Key" $runs4

